# ValueSet-KLConditionCodesHomeCare: publish version 1.1.0
# (commit: "Added 1.1.0 of term")
#
# The "Metadata" sheet holds a Property/Value table; bump the recorded
# Version and refresh the Date to match the new release.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
